# Generate Report for Archive
#
# 1. Update the report status text: the handoff has moved on, cells that
#    showed "Ready for handoff" now read "In Translation" (Overview summary
#    columns + the per-language Status columns all share this text).
# 2. Narrow the per-language "Status" column (and the Overview's matching
#    zh-cn/de-de summary columns) now that the new status text is shorter.

$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" --------------
# This single piece of text is shared by every sheet (Overview!E:F and the
# "Status" column on each language sheet), so sweep all worksheets.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Narrower Status-related columns ------------------------------------
# Target stored column width is ~13.41 character-units; Excel's ColumnWidth
# property is specified in whole characters, so use the character count that
# produces that stored width (12.5 characters).
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)
